$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 29
$wsExpo.Range("F4").Value = 210
$wsExpo.Range("F5").Value = 2592
$wsExpo.Range("F6").Value = 1859
$wsExpo.Range("F7").Value = 355
$wsExpo.Range("F8").Value = 110
$wsExpo.Range("F9").Value = 910

# Sheet "全部类型" (all types) - same underlying events, update column F accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 29
$wsAll.Range("F4").Value = 210
$wsAll.Range("F5").Value = 2592
$wsAll.Range("F6").Value = 1859
$wsAll.Range("F7").Value = 355
$wsAll.Range("F9").Value = 110
$wsAll.Range("F10").Value = 910
